# Generate Report for Handoff
# Updates the localization-status report with a new handoff/handback run:
#   - old GUID 0992185b-06cd-4604-b468-9574fe80adf6 -> new GUID fdfe4cc4-e100-4b84-b26c-6f430fac2abf
#   - old xliff hash c2e02d8fbcb85dba8775df5330a74bae91f98c53 -> new hash 0e677071c9f3a58be5de4ce09e103dfda5b1860f
#   - refreshed timestamps
#   - minor column-A width tweak (auto-fit) on all three sheets

$wb = $excel.ActiveWorkbook

$oldGuid = "0992185b-06cd-4604-b468-9574fe80adf6"
$newGuid = "fdfe4cc4-e100-4b84-b26c-6f430fac2abf"
$oldHash = "c2e02d8fbcb85dba8775df5330a74bae91f98c53"
$newHash = "0e677071c9f3a58be5de4ce09e103dfda5b1860f"

$oldMdName   = "$oldGuid.md"
$newMdName   = "$newGuid.md"
$oldMdPath   = "e2e\$oldGuid.md"
$newMdPath   = "e2e\$newGuid.md"

$oldZhXlf = "$oldGuid.$oldHash.zh-cn.xlf"
$newZhXlf = "$newGuid.$newHash.zh-cn.xlf"
$oldDeXlf = "$oldGuid.$oldHash.de-de.xlf"
$newDeXlf = "$newGuid.$newHash.de-de.xlf"

$oldHoDate  = "2016-08-30 04:57:14"
$newHoDate  = "2016-08-30 04:57:31"
$oldZhDate  = "2016-08-30 04:57:10"
$newZhDate  = "2016-08-30 04:57:26"

$newColWidth = 38.65   # produces a stored column width close to the refreshed auto-fit value

# The hyperlink target itself is not touched by this change (only the visible
# display text is refreshed) - it still points at the original blob URL.
$hyperlinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fe73b31bb65e5eb92e807e84c72402f3a577b1ab/e2e/$oldGuid.md"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("B2").Value = $newMdPath
$wsOverview.Range("G2").Value = $newHoDate

# Refresh the hyperlink display text on B2 (delete + re-add keeps the same
# target URL and relationship id, just refreshes the display text)
$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "$hyperlinkUrl", [Type]::Missing, [Type]::Missing, $newMdPath)

$wsOverview.Columns.Item(1).ColumnWidth = $newColWidth

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = $newMdName
$wsZhCn.Range("G2").Value = $newZhXlf
$wsZhCn.Range("H2").Value = $newZhDate

$wsZhCn.Range("A2").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "$hyperlinkUrl", [Type]::Missing, [Type]::Missing, $newMdName)

$wsZhCn.Columns.Item(1).ColumnWidth = $newColWidth

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = $newMdName
$wsDeDe.Range("G2").Value = $newDeXlf
$wsDeDe.Range("H2").Value = $newHoDate

$wsDeDe.Range("A2").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "$hyperlinkUrl", [Type]::Missing, [Type]::Missing, $newMdName)

$wsDeDe.Columns.Item(1).ColumnWidth = $newColWidth
